$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in worksheet order) that receive updated values for rows 2-25.
# Column G (value 1, slack-ish ref) and column H (unused/out-of-service bus)
# are unchanged.
$cols = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

# New vm_pu values per row (2..25), one array per row in the $cols order above.
$data = @{
    2  = @(1.02, 1.024555560652727, 1.032620573605937, 1.028181374237084, 1.039591645800745, 1.03047225865665, 1.029729517441904, 1.035425278111687, 1.030998938122216, 1.042376397342437, 1.013932003627131)
    3  = @(1.02, 1.025474795860711, 1.033309687536604, 1.029047846498326, 1.040474847788155, 1.030592728708559, 1.030287865905521, 1.035923522432904, 1.031673133753845, 1.043069657092351, 1.014119363382004)
    4  = @(1.02, 1.026070190676186, 1.033755895769136, 1.029609455629092, 1.041047153597179, 1.030669408074468, 1.0306491357337, 1.036245548855044, 1.032109689514306, 1.043518421212364, 1.014240522086863)
    5  = @(1.02, 1.026320634301036, 1.033943553252814, 1.029845780515499, 1.041287944453316, 1.030701338774864, 1.030801007948753, 1.036380838775689, 1.032293289590621, 1.043707123029228, 1.01429143868924)
    6  = @(1.02, 1.026362693042725, 1.033975065919389, 1.029885473619576, 1.041328385596682, 1.030706682159566, 1.030826507600928, 1.036403549237583, 1.032324121039214, 1.043738809280466, 1.0142999867127)
    7  = @(1.02, 1.026073536568899, 1.03375840297849, 1.029612612533902, 1.041050370298134, 1.030669835934426, 1.03065116508111, 1.03624735696065, 1.032112142506027, 1.043520942494436, 1.014241202510284)
    8  = @(1.02, 1.024866098009586, 1.03285339869223, 1.028474006305119, 1.039889958300007, 1.03051323526358, 1.029918217195763, 1.035593738217452, 1.031226721561176, 1.042610649875813, 1.013995338130202)
    9  = @(1.02, 1.022742993163858, 1.031261074578657, 1.026474930136489, 1.037851478733765, 1.030227562764897, 1.028626571882705, 1.034439189233432, 1.029668896241115, 1.041008028555446, 1.013561533415449)
    10 = @(1.02, 1.021330721758611, 1.030201238663796, 1.025147205377939, 1.036496834862475, 1.030030612908559, 1.027765473406548, 1.03366768427656, 1.028632033185811, 1.039940659243153, 1.013271978560796)
    11 = @(1.02, 1.020719949621527, 1.029742744118076, 1.024573488709939, 1.035911309388362, 1.029943796703766, 1.027392621225274, 1.03333320042613, 1.028183475638337, 1.039478741765706, 1.013146519913797)
    12 = @(1.02, 1.020493195540983, 1.029572503717037, 1.024360566104577, 1.035693977667355, 1.029911319158065, 1.027254129752064, 1.03320889671982, 1.028016924366157, 1.039307205564184, 1.013099907458586)
    13 = @(1.02, 1.020541829889685, 1.029609017899342, 1.02440623049767, 1.035740588860272, 1.029918296107145, 1.027283836516363, 1.033235563051753, 1.028052647310069, 1.039343998806044, 1.013109906503465)
    14 = @(1.02, 1.020701203728868, 1.029728670660614, 1.024555884753687, 1.035893341437556, 1.029941116791244, 1.027381173421258, 1.033322926691122, 1.028169707166503, 1.039464561695256, 1.013142667146476)
    15 = @(1.02, 1.020799414335202, 1.029802401321274, 1.02464811574985, 1.035987478394617, 1.02995514688887, 1.02744114625421, 1.033376746196402, 1.028241839988713, 1.039538849880568, 1.013162850544723)
    16 = @(1.02, 1.021371272604127, 1.030231676449847, 1.025185306408217, 1.036535716399647, 1.030036342308866, 1.027790218660469, 1.033689874178626, 1.028661811236756, 1.039971320790823, 1.013280303198705)
    17 = @(1.02, 1.021730185947906, 1.030501063331217, 1.025522593413183, 1.03687989233933, 1.030086863322556, 1.028009185750261, 1.033886180231622, 1.028925359005553, 1.040242669080398, 1.013353957157585)
    18 = @(1.02, 1.021939606386619, 1.030658232694687, 1.025719442601812, 1.037080744989219, 1.030116183172182, 1.02813690636032, 1.034000641908274, 1.029079121533317, 1.040400966986491, 1.013396910568856)
    19 = @(1.02, 1.022011025583981, 1.030711830216533, 1.025786582614179, 1.037149247574516, 1.03012615532818, 1.028180455893273, 1.034039663507466, 1.029131557256729, 1.040454946686787, 1.013411555253973)
    20 = @(1.02, 1.021691670450575, 1.030472156462289, 1.025486393764034, 1.036842955096107, 1.030081458221136, 1.027985692562274, 1.033865122612386, 1.028897078722062, 1.04021355337698, 1.013346055577977)
    21 = @(1.02, 1.02065426895284, 1.029693434100809, 1.024511810288337, 1.035848355222104, 1.029934403017045, 1.02735251005332, 1.033297201965546, 1.028135234208216, 1.039429057799774, 1.013133020268543)
    22 = @(1.02, 1.020002672914016, 1.029204196644609, 1.023900100715944, 1.03522392870269, 1.029840612326034, 1.026954417650447, 1.032939772400254, 1.027656597274337, 1.038936049394495, 1.012999010023417)
    23 = @(1.02, 1.020348033424736, 1.029463514491192, 1.024224279534162, 1.035554861467329, 1.029890458543089, 1.027165452231702, 1.033129285869531, 1.027910296766715, 1.039197379731544, 1.013070057554293)
    24 = @(1.02, 1.021709073722455, 1.030485218106953, 1.02550275047093, 1.036859645135569, 1.030083901011789, 1.027996308120274, 1.033874637768817, 1.028909857242715, 1.040226709432863, 1.013349625985509)
    25 = @(1.02, 1.023291319685245, 1.031672433374474, 1.026990865765381, 1.038377716264183, 1.030302564462645, 1.028960498123856, 1.034737991178851, 1.030071339465231, 1.041422166493374, 1.013673746152139)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
